$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting existing rows 74-195 down to 75-196
$ws.Rows.Item(74).Insert(-4121)  # xlShiftDown = -4121

# Populate the newly inserted row 74 with the new record's data
$ws.Cells.Item(74, 1).Value = 8
$ws.Cells.Item(74, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(74, 3).Value = "Coquimbo"
$ws.Cells.Item(74, 4).Value = 45100
$ws.Cells.Item(74, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(74, 5).Value = 4
$ws.Cells.Item(74, 6).Value = 100112052
$ws.Cells.Item(74, 7).Value = "Albahaca"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 800
$ws.Cells.Item(74, 11).Value = 3000
$ws.Cells.Item(74, 12).Value = 3500
$ws.Cells.Item(74, 13).Value = 3250
$ws.Cells.Item(74, 14).Value = "$/paquete"
$ws.Cells.Item(74, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(74, 16).Value = 3250
$ws.Cells.Item(74, 17).Value = 1
$ws.Cells.Item(74, 18).Value = "Hortaliza"
